$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 71; this shifts existing
# rows 71-182 down to 73-184 (matching dimension growing to A1:R184).
$ws.Rows("71:72").Insert()

# New row 71
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44792
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = 100112026
$ws.Range("G71").Value = "Haba"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 100
$ws.Range("K71").Value = 14000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = 14500
$ws.Range("N71").Value = "$/saco 25 kilos"
$ws.Range("O71").Value = "Provincia de Limarí"
$ws.Range("P71").Value = 580
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"

# New row 72
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44792
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 100112026
$ws.Range("G72").Value = "Haba"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Segunda"
$ws.Range("J72").Value = 25
$ws.Range("K72").Value = 12000
$ws.Range("L72").Value = 12000
$ws.Range("M72").Value = 12000
$ws.Range("N72").Value = "$/saco 25 kilos"
$ws.Range("O72").Value = "Provincia de Limarí"
$ws.Range("P72").Value = 480
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
